$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = "kWVo-vt4JQka9F6c5qtT7NYkj_A="
}

for ($r = 10; $r -le 18; $r++) {
    $ws.Cells.Item($r, 4).Value = "i32JKUsm7lIJ7ceMpSkZHlKq9cE="
}
